$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated daily case counts (column C) that cascade into the running-total
# formulas already present in column B.
$ws.Range("C372").Value = 60
$ws.Range("C374").Value = 92
$ws.Range("C376").Value = 29
$ws.Range("C377").Value = 82

# Row 378 (2021-03-09) newly filled in with real data; it previously held
# only the IF(TODAY()>...) placeholder formulas yielding blank strings.
$ws.Range("C378").Value = 13
$ws.Range("E378").Value = 9
$ws.Range("F378").Value = 8
$ws.Range("G378").Value = 38

# L378/M378 are formatted as Text (numFmtId 49), so a plain .Value = 0
# would be stored as the text string "0" instead of the number 0. Flip the
# number format to a numeric one long enough to write the value, then
# restore the original Text format by copying it back from a same-styled
# cell (avoids leaving stray/duplicate style records behind).
$ws.Range("D378").Copy()
$ws.Range("L378").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("L378").Value = 0
$ws.Range("L3").Copy()
$ws.Range("L378").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("F3").Copy()
$ws.Range("M378").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("M378").Value = 0
$ws.Range("M3").Copy()
$ws.Range("M378").PasteSpecial(-4122) # xlPasteFormats
